# "unify the conception of DataNode, DataTable, Entity."
# Rename the two sheets to match the unified naming scheme:
#   Property1 -> DataNode
#   Record    -> DataTable
# and make the DataTable sheet the active tab (it was the second sheet that
# ends up selected after the edit).
$wb = $excel.ActiveWorkbook

$wsNode  = $wb.Worksheets.Item("Property1")
$wsTable = $wb.Worksheets.Item("Record")

$wsNode.Name  = "DataNode"
$wsTable.Name = "DataTable"

$wsTable.Activate()
